$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Save", copying the formatting already used by G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for H2:H20 - default to 0, except H14 which is 1
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
